$d = $word.ActiveDocument

# Locate the paragraph that ends with the "LOQ4095..." requisito line; the
# three paragraphs that directly follow it (a blank paragraph, a
# page-break paragraph, and the "© 2020 ..." footer paragraph) are being
# removed, while the trailing blank + page-break paragraphs at the very
# end of the document must be kept untouched.
$anchor = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "LOQ4095: Qu.mica Geral Experimental \(Requisito\)") {
        $anchor = $i
    }
}

if ($anchor -eq $null) {
    throw "Could not find the LOQ4095 requisito paragraph"
}

$startPara = $d.Paragraphs.Item($anchor + 1)
$endPara = $d.Paragraphs.Item($anchor + 3)

$rng = $d.Range($startPara.Range.Start, $endPara.Range.End)
$rng.Delete()

Write-Output "Removed paragraphs $($anchor + 1)..$($anchor + 3); new paragraph count: $($d.Paragraphs.Count)"
